$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.330.44'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '2.270.26'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '499.62'
$ws.Range("E5").Value = '  +1.53%  '

$ws.Range("D6").Value = '129.23'
$ws.Range("E6").Value = '  +1.76%  '

$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +0.39%  '

$ws.Range("E10").Value = '  +0.97%  '

$ws.Range("D11").Value = '0.337'
$ws.Range("E11").Value = '  +3.98%  '

$ws.Range("D12").Value = '4.90'
$ws.Range("E12").Value = '  +5.43%  '

$ws.Range("D13").Value = '23.22'
$ws.Range("E13").Value = '  +7.11%  '

$ws.Range("D14").Value = '2.671.73'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = '54.292.73'
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("E16").Value = '  +1.10%  '

$ws.Range("D17").Value = '2.279.71'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").Value = '10.29'
$ws.Range("E18").Value = '  +3.03%  '

$ws.Range("E19").Value = '  +1.99%  '

$ws.Range("D20").Value = '303.87'
$ws.Range("E20").Value = '  +2.00%  '

$ws.Range("D21").Value = '6.31'
$ws.Range("E21").Value = '  -1.19%  '

$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = '60.36'
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").Value = '0.996'
$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("E25").Value = '  +1.39%  '

$ws.Range("D26").Value = '7.36'
$ws.Range("E26").Value = '  +5.05%  '

$ws.Range("D27").Value = '175.54'
$ws.Range("E27").Value = '  +5.46%  '

$ws.Range("E28").Value = '  +3.40%  '

$ws.Range("D29").Value = '6.01'
$ws.Range("E29").Value = '  +3.32%  '

$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("E31").Value = '  +2.53%  '

$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("D33").Value = '17.80'
$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("E35").Value = '  +6.15%  '

$ws.Range("E36").Value = '  +2.20%  '

$ws.Range("D37").Value = '3.74'
$ws.Range("E37").Value = '  +2.18%  '

$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("D39").Value = '1.41'
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  +1.36%  '

$ws.Range("E41").Value = '  +2.35%  '

$ws.Range("D42").Value = '124.97'
$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = '0.0492'
$ws.Range("E43").Value = '  +2.24%  '

$ws.Range("E44").Value = '  +1.42%  '

$ws.Range("D45").Value = '245.53'
$ws.Range("E45").Value = '  +4.26%  '

$ws.Range("E46").Value = '  +1.24%  '

$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("E48").Value = '  +1.95%  '

$ws.Range("E49").Value = '  +0.85%  '

$ws.Range("E50").Value = '  +1.18%  '

$ws.Range("D51").Value = '1.53'
$ws.Range("E51").Value = '  +3.44%  '
